$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the data-row template placeholders from upper-case field names
# (row.ORDER_ID, row.CITY_NAME, row.ITEM_NAME, row.ORDER_DATE, row.VOLUME)
# to lower-case field names (row.order_id, row.city_name, row.item_name,
# row.order_date, row.volume).
$ws.Range("A10").Value = '${row.order_id}'
$ws.Range("B10").Value = '${row.city_name}'
$ws.Range("C10").Value = '${row.item_name}'
$ws.Range("D10").Value = '${row.order_date}'
$ws.Range("E10").Value = '${row.volume}'

# Move the active selection to F11, matching the saved workbook state.
$ws.Range("F11").Select()
